# Remove all comments from every slide in the presentation.
$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    $count = $s.Comments.Count
    for ($j = $count; $j -ge 1; $j--) {
        $s.Comments.Item($j).Delete()
    }
}
